# Add 2022-Q4 data:
#  - Insert a new worksheet "2022-Q4" right after "总计" (before the current "2022-Q3" sheet),
#    built as a duplicate of "2022-Q3" (so it inherits the same column layout/styling),
#    then overwrite its cells with the Q4 fund-holding data.
#  - Update the "总计" (summary) sheet: the former top row becomes the 2022-Q4 entry and a
#    fresh row is inserted carrying the old 2022-Q3 totals (the 2022-Q2 row shifts down).

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet before the existing "2022-Q3" sheet
#    by duplicating "2022-Q3" (keeps header/row styling identical) and
#    renaming the copy.
# ------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($q3Sheet)
$q4Sheet = $wb.Worksheets.Item("2022-Q3 (2)")
$q4Sheet.Name = "2022-Q4"

# Columns B and D-G are stored as plain text in the source data (fund codes must
# keep their leading zeros, and the numeric-looking columns must keep trailing
# zeros / exact decimal formatting) - prefix with an apostrophe so Excel stores
# them as text instead of coercing to a number. Column C (fund name) and column A
# (0-based index) / H (rank) keep their original types.
$q4Data = @(
    @("010490", "鹏华高质量增长混合A",   "13.18", "94.62", "4.48", "0.5905", 8),
    @("009023", "鹏华稳健回报混合A",     "4.33",  "94.26", "5.49", "0.2377", 7),
    @("014541", "华安新能源主题混合A",   "1.09",  "90.05", "3.10", "0.0338", 8),
    @("010491", "鹏华高质量增长混合C",   "0.45",  "94.62", "4.48", "0.0202", 8),
    @("014542", "华安新能源主题混合C",   "0.10",  "90.05", "3.10", "0.0031", 8),
    @("017511", "鹏华稳健回报混合C",     "0.05",  "94.26", "5.49", "0.0027", 7)
)

for ($r = 0; $r -lt $q4Data.Length; $r++) {
    $row = $r + 2
    $rec = $q4Data[$r]
    $q4Sheet.Cells.Item($row, 1).Value = $r
    $q4Sheet.Cells.Item($row, 2).Value = "'" + $rec[0]
    $q4Sheet.Cells.Item($row, 3).Value = $rec[1]
    $q4Sheet.Cells.Item($row, 4).Value = "'" + $rec[2]
    $q4Sheet.Cells.Item($row, 5).Value = "'" + $rec[3]
    $q4Sheet.Cells.Item($row, 6).Value = "'" + $rec[4]
    $q4Sheet.Cells.Item($row, 7).Value = "'" + $rec[5]
    $q4Sheet.Cells.Item($row, 8).Value = $rec[6]
}

# ------------------------------------------------------------------
# 2. Update the "总计" summary sheet
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Make room for the new top row by pushing the old rows down one slot
$totalSheet.Rows.Item(3).Insert()

# Row 2 now becomes the 2022-Q4 summary entry
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 0.89

# Row 3 carries what used to be the 2022-Q3 summary entry
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 6
$totalSheet.Range("D3").Value = 1.31

# Row 4 (previously row 3, shifted down) keeps the 2022-Q2 values, just renumber A4
$totalSheet.Range("A4").Value = 2

# ------------------------------------------------------------------
# 3. Restore the originally-active tab ("2022-Q2") - adding/renaming sheets
#    above shifted the active tab onto the new "2022-Q4" sheet.
# ------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q2").Activate()
